$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 157.6
$ws.Range("I38").Value = 157.6
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 472.8
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -100.8
$ws.Range("N38").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 755.6667
$ws.Range("I58").Value = 125
$ws.Range("J58").Value = 2017
$ws.Range("K58").Value = 375
$ws.Range("L58").Value = 6051
$ws.Range("M58").Value = -225
$ws.Range("N58").Value = -6351

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3333.3333
$ws.Range("I74").Value = 2666.6667
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 2666.6667
$ws.Range("L74").Value = 4000
$ws.Range("M74").Value = -1730.6667
$ws.Range("N74").Value = -5872

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3333.3333
$ws.Range("I77").Value = 2666.6667
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 13333.3335
$ws.Range("L77").Value = 20000
$ws.Range("M77").Value = -8653.333500000001
$ws.Range("N77").Value = -29360

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 39800
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 39800
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 39800
$ws.Range("N87").Value = -42296

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 39800
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 39800
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 119400
$ws.Range("N90").Value = -131880

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1444121.8
$ws.Range("I132").Value = 2925.9333
$ws.Range("J132").Value = 12253090
$ws.Range("K132").Value = 8777.7999
$ws.Range("L132").Value = 36759270
$ws.Range("M132").Value = -6247.7999
$ws.Range("N132").Value = -36764330

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 7502187.5
$ws.Range("I137").Value = 16668356
$ws.Range("J137").Value = 2595.5
$ws.Range("K137").Value = 50005068
$ws.Range("L137").Value = 7786.5
$ws.Range("M137").Value = -50002518
$ws.Range("N137").Value = -12886.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6311718.5
$ws.Range("I32").Value = 6826168.5
$ws.Range("J32").Value = 9707
$ws.Range("K32").Value = 6826168.5
$ws.Range("L32").Value = 9707
$ws.Range("M32").Value = -6825881.5
$ws.Range("N32").Value = -10281

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I63").Value = 2399.0833
$ws.Range("J63").Value = 4506
$ws.Range("K63").Value = 2399.0833
$ws.Range("L63").Value = 4506
$ws.Range("M63").Value = -1713.0833
$ws.Range("N63").Value = -5878

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I66").Value = 2399.0833
$ws.Range("J66").Value = 4506
$ws.Range("K66").Value = 11995.4165
$ws.Range("L66").Value = 22530
$ws.Range("M66").Value = -8563.416499999999
$ws.Range("N66").Value = -29394

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 28224.334
$ws.Range("I132").Value = 21308.84
$ws.Range("J132").Value = 42055.32
$ws.Range("K132").Value = 63926.52
$ws.Range("L132").Value = 126165.96
$ws.Range("M132").Value = -61396.52
$ws.Range("N132").Value = -131225.96

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 10147.621
$ws.Range("I86").Value = 21958
$ws.Range("J86").Value = 2930.1667
$ws.Range("K86").Value = 21958
$ws.Range("L86").Value = 2930.1667
$ws.Range("M86").Value = -20835
$ws.Range("N86").Value = -5176.1667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 10147.621
$ws.Range("I89").Value = 21958
$ws.Range("J89").Value = 2930.1667
$ws.Range("K89").Value = 109790
$ws.Range("L89").Value = 14650.8335
$ws.Range("M89").Value = -104174
$ws.Range("N89").Value = -25882.8335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 25944.166
$ws.Range("I31").Value = 33858.59
$ws.Range("J31").Value = 6723.4287
$ws.Range("K31").Value = 33858.59
$ws.Range("L31").Value = 6723.4287
$ws.Range("M31").Value = -33563.59
$ws.Range("N31").Value = -7313.4287

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 25944.166
$ws.Range("I34").Value = 33858.59
$ws.Range("J34").Value = 6723.4287
$ws.Range("K34").Value = 33858.59
$ws.Range("L34").Value = 6723.4287
$ws.Range("M34").Value = -33656.59
$ws.Range("N34").Value = -7127.4287

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H81").Value = 42000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 42000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 42000
$ws.Range("N81").Value = -43996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H84").Value = 42000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 42000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 126000
$ws.Range("N84").Value = -135984

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 995.3333
$ws.Range("I17").Value = 746
$ws.Range("J17").Value = 1120
$ws.Range("K17").Value = 2238
$ws.Range("L17").Value = 3360
$ws.Range("M17").Value = -2069
$ws.Range("N17").Value = -3698

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 166667380
$ws.Range("I112").Value = 1397
$ws.Range("J112").Value = 333333340
$ws.Range("K112").Value = 4191
$ws.Range("L112").Value = 1000000020
$ws.Range("M112").Value = -3083
$ws.Range("N112").Value = -1000002236

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 290.33334
$ws.Range("I2").Value = 290.33334
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 290.33334
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -177.33334
$ws.Range("N2").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 3000
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 3000
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -3346

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 3000
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 3000
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = -3210

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 2628.8333
$ws.Range("I31").Value = 2628.8333
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2628.8333
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -2336.8333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H37").Value = 2628.8333
$ws.Range("I37").Value = 2628.8333
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 2628.8333
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -2351.8333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 1000
$ws.Range("I43").Value = 1000
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 1000
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -849
$ws.Range("N43").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2697.1765
$ws.Range("I80").Value = 1801.125
$ws.Range("J80").Value = 3493.6667
$ws.Range("K80").Value = 1801.125
$ws.Range("L80").Value = 3493.6667
$ws.Range("M80").Value = -803.125
$ws.Range("N80").Value = -5489.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2697.1765
$ws.Range("I83").Value = 1801.125
$ws.Range("J83").Value = 3493.6667
$ws.Range("K83").Value = 9005.625
$ws.Range("L83").Value = 17468.3335
$ws.Range("M83").Value = -4013.625
$ws.Range("N83").Value = -27452.3335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 34922.74
$ws.Range("I132").Value = 31387.646
$ws.Range("J132").Value = 39374.332
$ws.Range("K132").Value = 94162.93799999999
$ws.Range("L132").Value = 118122.996
$ws.Range("M132").Value = -91632.93799999999
$ws.Range("N132").Value = -123182.996

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2610
$ws.Range("I7").Value = 2540
$ws.Range("J7").Value = 2680
$ws.Range("K7").Value = 2540
$ws.Range("L7").Value = 2680
$ws.Range("M7").Value = -2428
$ws.Range("N7").Value = -2904

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1501
$ws.Range("I46").Value = 375.25
$ws.Range("J46").Value = 3002
$ws.Range("K46").Value = 375.25
$ws.Range("L46").Value = 3002
$ws.Range("M46").Value = -187.25
$ws.Range("N46").Value = -3378

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2610
$ws.Range("I126").Value = 2540
$ws.Range("J126").Value = 2680
$ws.Range("K126").Value = 7620
$ws.Range("L126").Value = 8040
$ws.Range("M126").Value = -5150
$ws.Range("N126").Value = -12980

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 63064.938
$ws.Range("I136").Value = 29518.111
$ws.Range("J136").Value = 172854.55
$ws.Range("K136").Value = 88554.333
$ws.Range("L136").Value = 518563.65
$ws.Range("M136").Value = -86004.333
$ws.Range("N136").Value = -523663.65

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4699.857
$ws.Range("I62").Value = 3450
$ws.Range("J62").Value = 5199.8
$ws.Range("K62").Value = 3450
$ws.Range("L62").Value = 5199.8
$ws.Range("M62").Value = -2826
$ws.Range("N62").Value = -6447.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 4699.857
$ws.Range("I65").Value = 3450
$ws.Range("J65").Value = 5199.8
$ws.Range("K65").Value = 17250
$ws.Range("L65").Value = 25999
$ws.Range("M65").Value = -14130
$ws.Range("N65").Value = -32239

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 34650.625
$ws.Range("I82").Value = 31001
$ws.Range("J82").Value = 35172
$ws.Range("K82").Value = 31001
$ws.Range("L82").Value = 35172
$ws.Range("M82").Value = -30618
$ws.Range("N82").Value = -35938

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H85").Value = 34650.625
$ws.Range("I85").Value = 31001
$ws.Range("J85").Value = 35172
$ws.Range("K85").Value = 31001
$ws.Range("L85").Value = 35172
$ws.Range("M85").Value = -29675
$ws.Range("N85").Value = -37824

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1473.8889
$ws.Range("I113").Value = 1033.125
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 3099.375
$ws.Range("L113").Value = 15000
$ws.Range("M113").Value = -929.375
$ws.Range("N113").Value = -19340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 62528.848
$ws.Range("I132").Value = 45241.523
$ws.Range("J132").Value = 102289.7
$ws.Range("K132").Value = 135724.569
$ws.Range("L132").Value = 306869.1
$ws.Range("M132").Value = -133194.569
$ws.Range("N132").Value = -311929.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 46427.637
$ws.Range("I136").Value = 31240.94
$ws.Range("J136").Value = 91987.73
$ws.Range("K136").Value = 93722.81999999999
$ws.Range("L136").Value = 275963.19
$ws.Range("M136").Value = -91172.81999999999
$ws.Range("N136").Value = -281063.19
